$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Strip the old border/number-format style off E:H (cols lose the old numFmt 176 + border),
# using a blank helper cell cut onto each target cell; Cut only carries over formatting
# from the (empty) source, leaving the destination value untouched.
$blank = $ws.Range("ZZ9999")
$blank.Cut($ws.Range("E1:E13"))
$blank.Cut($ws.Range("F1:F13"))
$blank.Cut($ws.Range("G1:G13"))
$blank.Cut($ws.Range("H1:H13"))

# --- Write the corrected E/F/G/H values (F/G/H were rotated) and the new I (serving size) column
$ws.Range("E1").Value = 0.7376
$ws.Range("F1").Value = 0.07940000000000001
$ws.Range("G1").Value = 0.067
$ws.Range("H1").Value = 0.017
$ws.Range("I1").Value = 500
$ws.Range("E2").Value = 1.489025
$ws.Range("F2").Value = 0.11225
$ws.Range("G2").Value = 0.11474999999999999
$ws.Range("H2").Value = 0.0645
$ws.Range("I2").Value = 400
$ws.Range("E3").Value = 1.8615666666666668
$ws.Range("F3").Value = 0.077
$ws.Range("G3").Value = 0.15166666666666667
$ws.Range("H3").Value = 0.10533333333333333
$ws.Range("I3").Value = 300
$ws.Range("E4").Value = 2.5246142857142857
$ws.Range("F4").Value = 0.1906857142857143
$ws.Range("G4").Value = 0.16517142857142855
$ws.Range("H4").Value = 0.12235714285714284
$ws.Range("I4").Value = 70
$ws.Range("E5").Value = 1.84
$ws.Range("F5").Value = 0.311
$ws.Range("G5").Value = 0.031
$ws.Range("H5").Value = 0.052000000000000005
$ws.Range("I5").Value = 100
$ws.Range("E6").Value = 1.588
$ws.Range("F6").Value = 0.008400000000000001
$ws.Range("G6").Value = 0.2444
$ws.Range("H6").Value = 0.066
$ws.Range("I6").Value = 250
$ws.Range("E7").Value = 2.4032
$ws.Range("F7").Value = 0.081
$ws.Range("G7").Value = 0.195
$ws.Range("H7").Value = 0.14400000000000002
$ws.Range("I7").Value = 100
$ws.Range("E8").Value = 1.952
$ws.Range("F8").Value = 0.0008
$ws.Range("G8").Value = 0.2184
$ws.Range("H8").Value = 0.1232
$ws.Range("I8").Value = 250
$ws.Range("E9").Value = 2.635266666666667
$ws.Range("F9").Value = 0.052
$ws.Range("G9").Value = 0.22933333333333333
$ws.Range("H9").Value = 0.16799999999999998
$ws.Range("I9").Value = 150
$ws.Range("E10").Value = 2.9654000000000003
$ws.Range("F10").Value = 0.11599999999999999
$ws.Range("G10").Value = 0.26
$ws.Range("H10").Value = 0.162
$ws.Range("I10").Value = 100
$ws.Range("E11").Value = 3.2971666666666666
$ws.Range("F11").Value = 0.08733333333333333
$ws.Range("G11").Value = 0.20033333333333334
$ws.Range("H11").Value = 0.23866666666666664
$ws.Range("I11").Value = 300
$ws.Range("E12").Value = 3.2054
$ws.Range("F12").Value = 0.09133333333333332
$ws.Range("G12").Value = 0.18066666666666667
$ws.Range("H12").Value = 0.2353333333333333
$ws.Range("I12").Value = 150
$ws.Range("E13").Value = 2.889013333333333
$ws.Range("F13").Value = 0.05869333333333333
$ws.Range("G13").Value = 0.20514000000000002
$ws.Range("H13").Value = 0.20374
$ws.Range("I13").Value = 150

# --- Give the new I column (rows 1-13) the thin border the old E:H columns used to carry
$ws.Range("I1:I13").Borders.LineStyle = 1

# --- Remove the now-orphaned last data row (old D000023 / 임연수구이 entry):
# clear A:D (and I, never populated) completely so they disappear, and blank out E:H
$blank.Cut($ws.Range("A14"))
$blank.Cut($ws.Range("B14"))
$blank.Cut($ws.Range("C14"))
$blank.Cut($ws.Range("D14"))
$ws.Range("A14:D14").ClearContents()
$ws.Range("E14:H14").ClearContents()

# --- Restore the selection the author left active
$ws.Range("G16").Select()
